$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ShipmentTracking values (column P). The tracking numbers are
# all-digit strings, so they must be forced to remain text (shared string)
# instead of being auto-converted to numbers by the COM layer. We do this
# by temporarily switching the cell to a text number format, assigning the
# value, and then resetting the cell style back to Normal so no stray
# style index is left referenced on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("P2") "320018631291"
Set-TextValue $ws.Range("P3") "320018638712"
Set-TextValue $ws.Range("P5") "320018639318"
